# "Generate Report for Archive"
#
# The localization status changes from "Ready for handoff" to
# "In Translation" for the dd1ade58-...md file. That string shows up in
# three places:
#   - Overview!E2  (zh-cn status)
#   - Overview!F2  (de-de status)
#   - zh-cn!C2     (Status column)
#   - de-de!C2     (Status column)
#
# Because "Ready for handoff" is a shared string, retyping it everywhere it
# appears collapses back down to a single shared-string entry on save (same
# as the original file), and since "In Translation" is shorter than
# "Ready for handoff" the Status/locale columns that held it get narrower to
# fit the new text (columns E & F on Overview, column C on zh-cn/de-de).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (E) / de-de (F) status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns("E:F").ColumnWidth = 12.5

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns("C:C").ColumnWidth = 12.5

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns("C:C").ColumnWidth = 12.5

Write-Output "Updated status to 'In Translation' on Overview!E2:F2, zh-cn!C2, de-de!C2 and resized their columns."
